$d = $word.ActiveDocument

# 1. Merge "Professor Vajda Weekly Meeting #" + "8" -> single run text stays the same when read, but
#    we want the underlying OOXML to merge. Using Find/Replace on "#8" with formatted text keeps it simple.
$d.Content.Find.Execute("Professor Vajda Weekly Meeting #8", $true, $false, $false, $false, $false, $true, 1, $false, "Professor Vajda Weekly Meeting #8", 2)

# 2. Merge "March 7th, 2023 / 10:00 AM / Zoom"
$d.Content.Find.Execute("March 7th, 2023 / 10:00 AM / Zoom", $true, $false, $false, $false, $false, $true, 1, $false, "March 7th, 2023 / 10:00 AM / Zoom", 2)

# 3. Merge Attendees paragraph runs
$d.Content.Find.Execute("Professor Szilard Vajda, Presentation Spectators", $true, $false, $false, $false, $false, $true, 1, $false, "Professor Szilard Vajda, Presentation Spectators", 2)

# 4. Group Presentation (40 Minutes)
$d.Content.Find.Execute("Group Presentation (40 Minutes)", $true, $false, $false, $false, $false, $true, 1, $false, "Group Presentation (40 Minutes)", 2)

# 5. Questions and Answers (10 Minutes)
$d.Content.Find.Execute("Questions and Answers (10 Minutes)", $true, $false, $false, $false, $false, $true, 1, $false, "Questions and Answers (10 Minutes)", 2)

# 6. Peer review sentence update
$d.Content.Find.Execute("Should each turn in a peer review for Assignment 2.", $true, $false, $false, $false, $false, $true, 1, $false, "Should each turn in a peer review for Assignment 2 by Wednesday, March 8th.", 2)
